$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the unused "Not Used" header cell at I1 entirely (column I is no longer used).
$ws.Range("I1").Clear()

# Add new "Is Active?" header in K1, copying the bold/underline header formatting
# from the neighboring J1 header cell so the style matches the other headers.
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "Is Active?"

# Column K used to hold the (unused) wide "Not Used" column; now it holds short
# TRUE/FALSE values, so shrink it to a narrower width.
$ws.Columns.Item(11).ColumnWidth = 10.5

# Populate the Is Active? column with TRUE for every existing data row.
$ws.Range("K2").Value = $true
$ws.Range("K3").Value = $true
$ws.Range("K5").Value = $true
$ws.Range("K6").Value = $true

# Update the active selection to L1, as in the saved workbook.
[void]$ws.Range("L1").Select()
